$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.583.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.729.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2664"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06177"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07107"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.532"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.588.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006966"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.955.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.515"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.786"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.240"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.405"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.766"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.963"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08007"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.676"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04563"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6299"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8957"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.017"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.384"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.004"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.396"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3886"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1183"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05391"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.868"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.67%  "
